$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.101702451705933
$ws.Range("B1").Value = 2.092604637145996
$ws.Range("C1").Value = 9.370265960693359
$ws.Range("D1").Value = 2.406872510910034
$ws.Range("E1").Value = 1.292432427406311
